# Daily attendance processing - 2026-01-20 10:42:21
#
# The "Recorded By" column (G) on the "Session Analysis Results" sheet lists
# who/what recorded each attendance session, as a comma-separated list (e.g.
# "dnasr281@gmail.com, System"). Going forward "System" should be reported
# first in that list, so for every row where "System" (exact case) appears
# alongside other recorders, move it to the front while keeping the relative
# order of the remaining recorders unchanged. Note "system" (lower case) is
# a distinct recorder and must be left where it is.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# PowerShell's default -eq/-ne/-contains operators here are case-INsensitive,
# so "System" -eq "system" would be True. Do a case-sensitive compare via a
# regex match anchored with the (?-i) inline flag instead.
function Test-ExactMatch($text, $literal) {
    $escaped = [regex]::Escape($literal)
    return ($text -cmatch "(?-i)^$escaped`$")
}

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $usedRange.Rows.Count + $firstRow - 1

for ($r = 2; $r -le $lastRow; $r++) {

    # Column I ("Status") is populated for every data row; only rows that
    # were actually recorded have a "Recorded By" (column G) value at all,
    # so gate on it rather than touching an empty G cell directly.
    $status = $ws.Cells.Item($r, 9).Text
    if ($status -ne "Recorded") {
        continue
    }

    $cell = $ws.Cells.Item($r, 7)
    $value = $cell.Text

    if ([string]::IsNullOrEmpty($value)) {
        continue
    }

    $parts = $value -split ", "

    $hasSystem = $false
    $hasAdmin = $false
    foreach ($p in $parts) {
        if (Test-ExactMatch $p "System") { $hasSystem = $true }
        if (Test-ExactMatch $p "admin@admin.com") { $hasAdmin = $true }
    }

    if ($parts.Count -gt 1 -and $hasSystem -and -not $hasAdmin) {
        $rest = @()
        $removedOne = $false
        foreach ($p in $parts) {
            if ((-not $removedOne) -and (Test-ExactMatch $p "System")) {
                $removedOne = $true
                continue
            }
            $rest += $p
        }
        $newParts = @("System") + $rest
        $newValue = $newParts -join ", "
        $cell.Value = $newValue
    }
}
